# Update the "Ready for handoff" status text to "In Translation" everywhere
# it appears, and shrink the width of the status-related columns across the
# Overview, zh-cn, and de-de worksheets.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        if ("Ready for handoff" -eq $cell.Value2) {
            $cell.Value = "In Translation"
        }
    }
}

$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 12.5
